$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Igf2", "Igf1r", "ECs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 14.561928, 43.685784, 0.3501777048818433, 0.3501777048818433, 124.412152068528, 1119.709368616752, 0.008928349922480351, 0.008928349922480351),
    @("ECs", "Igf2", "Igf1r", "FAPs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 14.40015733333333, 43.200472, 0.3462875276490937, 0.3462875276490937, 123.0300386024018, 1107.270347421616, 0.008829163529085677, 0.008829163529085677),
    @("ECs", "Igf2", "Igf1r", "M2", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 2.886742333333333, 8.660226999999999, 0.06941888497676431, 0.06941888497676431, 24.66334307911177, 221.970087712006, 0.001769947337195831, 0.001769947337195831),
    @("ECs", "Igf2", "Igf1r", "sCs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 9.735567333333334, 29.206702, 0.2341158824922987, 0.2341158824922987, 83.17737071272845, 748.5963364145559, 0.005969165061513071, 0.005969165061513071),
    @("FAPs", "Igf2", "Igf1r", "ECs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 14.561928, 43.685784, 0.3501777048818433, 0.3501777048818433, 3965.61951489, 35690.57563401, 0.2845898740570963, 0.2845898740570963),
    @("FAPs", "Igf2", "Igf1r", "FAPs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 14.40015733333333, 43.200472, 0.3462875276490937, 0.3462875276490937, 3921.564846258889, 35294.08361633, 0.2814283219842665, 0.2814283219842665),
    @("FAPs", "Igf2", "Igf1r", "M2", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 2.886742333333333, 8.660226999999999, 0.06941888497676431, 0.06941888497676431, 786.1405255901389, 7075.264730311249, 0.05641681768228916, 0.05641681768228915),
    @("FAPs", "Igf2", "Igf1r", "sCs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 9.735567333333334, 29.206702, 0.2341158824922987, 0.2341158824922987, 2651.266769454723, 23861.4009250925, 0.1902662807608796, 0.1902662807608796),
    @("sCs", "Igf2", "Igf1r", "ECs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 14.561928, 43.685784, 0.3501777048818433, 0.3501777048818433, 38.68126193435999, 348.13135740924, 0.002775933349363301, 0.002775933349363301),
    @("sCs", "Igf2", "Igf1r", "FAPs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 14.40015733333333, 43.200472, 0.3462875276490937, 0.3462875276490937, 38.25154592899111, 344.26391336092, 0.002745095084777133, 0.002745095084777133),
    @("sCs", "Igf2", "Igf1r", "M2", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 2.886742333333333, 8.660226999999999, 0.06941888497676431, 0.06941888497676431, 7.668135451066109, 69.01321905959499, 0.0005502983062489273, 0.0005502983062489273),
    @("sCs", "Igf2", "Igf1r", "sCs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 9.735567333333334, 29.206702, 0.2341158824922987, 0.2341158824922987, 25.86086334860778, 232.74777013747, 0.001855886530655277, 0.001855886530655277),
    @("M2", "Igf2", "Igf1r", "ECs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 14.561928, 43.685784, 0.3501777048818433, 0.3501777048818433, 750.8406559272961, 6757.565903345664, 0.05388354755290335, 0.05388354755290333),
    @("M2", "Igf2", "Igf1r", "FAPs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 14.40015733333333, 43.200472, 0.3462875276490937, 0.3462875276490937, 742.4994532053903, 6682.495078848513, 0.05328494705096444, 0.05328494705096443),
    @("M2", "Igf2", "Igf1r", "M2", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 2.886742333333333, 8.660226999999999, 0.06941888497676431, 0.06941888497676431, 148.8459156681102, 1339.613241012992, 0.0106818216510304, 0.0106818216510304),
    @("M2", "Igf2", "Igf1r", "sCs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 9.735567333333334, 29.206702, 0.2341158824922987, 0.2341158824922987, 501.9843363038436, 4517.859026734593, 0.03602455013925072, 0.03602455013925072)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, 1 + $j).Value = $row[$j]
    }
}
